# Insert a new data row at row 26 (shifting existing rows 26-115 down to 27-116)
# and populate it with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(26).Insert()

$ws.Range("A26").Value2 = 7
$ws.Range("B26").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C26").Value2 = "Ñuble"
$ws.Range("D26").Value2 = 44459
$ws.Range("E26").Value2 = 16
$ws.Range("F26").Value2 = 100112017
$ws.Range("G26").Value2 = "Apio"
$ws.Range("H26").Value2 = "Americana (o)"
$ws.Range("I26").Value2 = "Primera"
$ws.Range("J26").Value2 = 100
$ws.Range("K26").Value2 = 8000
$ws.Range("L26").Value2 = 8500
$ws.Range("M26").Value2 = 8250
$ws.Range("N26").Value2 = "`$/docena de matas"
$ws.Range("O26").Value2 = "Provincia del Elquí"
$ws.Range("P26").Value2 = 1375
$ws.Range("Q26").Value2 = 6
$ws.Range("R26").Value2 = "Hortaliza"
